# Apply crypto price/volume updates (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.852.99'
$ws.Range('E2').Value = '  +1.38%  '
$ws.Range('D3').Value = '3.312.90'
$ws.Range('E3').Value = '  +6.24%  '
$ws.Range('D5').Value = "'601.44"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.39%  '
$ws.Range('D6').Value = "'142.79"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.86%  '
$ws.Range('D8').Value = '3.310.10'
$ws.Range('E8').Value = '  +6.28%  '
$ws.Range('D9').Value = "'0.521"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.45%  '
$ws.Range('E10').Value = '  +2.97%  '
$ws.Range('D11').Value = "'5.54"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.45%  '
$ws.Range('E12').Value = '  +4.03%  '
$ws.Range('E13').Value = '  +1.31%  '
$ws.Range('D14').Value = "'34.79"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.93%  '
$ws.Range('D15').Value = '3.859.11'
$ws.Range('E15').Value = '  +6.32%  '
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('D17').Value = '3.310.84'
$ws.Range('E17').Value = '  +6.44%  '
$ws.Range('D18').Value = '63.936.00'
$ws.Range('E18').Value = '  +1.48%  '
$ws.Range('D19').Value = "'6.91"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.66%  '
$ws.Range('D20').Value = "'481.07"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.71%  '
$ws.Range('D21').Value = "'14.22"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.95%  '
$ws.Range('E22').Value = '  +5.73%  '
$ws.Range('E23').Value = '  +4.60%  '
$ws.Range('D24').Value = "'84.97"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.20%  '
$ws.Range('D25').Value = "'13.49"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.39%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  +2.09%  '
$ws.Range('D28').Value = "'7.32"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.59%  '
$ws.Range('E30').Value = '  +3.60%  '
$ws.Range('E31').Value = '  +4.22%  '
$ws.Range('D32').Value = "'29.45"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +10.39%  '
$ws.Range('E33').Value = '  -1.37%  '
$ws.Range('E34').Value = '  +1.05%  '
$ws.Range('E35').Value = '  +2.46%  '
$ws.Range('D37').Value = "'52.94"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.82%  '
$ws.Range('D38').Value = '0.0₃0753'
$ws.Range('E38').Value = '  +7.96%  '
$ws.Range('E39').Value = '  +4.53%  '
$ws.Range('D40').Value = "'431.16"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.39%  '
$ws.Range('D41').Value = '3.046.89'
$ws.Range('E41').Value = '  +5.06%  '
$ws.Range('E42').Value = '  +2.74%  '
$ws.Range('E43').Value = '  +2.61%  '
$ws.Range('D45').Value = "'0.266"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.51%  '
$ws.Range('E46').Value = '  +4.56%  '
$ws.Range('D47').Value = "'26.43"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.83%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').Value = "'0.999"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D49').Value = "'35.91"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +14.67%  '
$ws.Range('E50').Value = '  +2.17%  '
$ws.Range('E51').Value = '  +3.28%  '
